# Updated versions based on review comments
#
# Lower-cases several labels/phrases on the "Meta-release journey" slide
# (slide 3) and tweaks the release-cadence wording, plus updates the
# title of the "Beyond Meta-release" slide (slide 5) to call out Fall24
# explicitly.

function Replace-Substring {
    param(
        $Shape,
        [string]$Old,
        [string]$New
    )
    $tr = $Shape.TextFrame.TextRange
    $full = $tr.Text
    $idx = $full.IndexOf($Old)
    if ($idx -lt 0) {
        Write-Host "NOT FOUND: '$Old' in shape '$($Shape.Name)' text '$full'"
        return
    }
    $sub = $tr.Characters($idx + 1, $Old.Length)
    $sub.Text = $New
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# Slide 3 - "Meta-release journey" slide
# ---------------------------------------------------------------
$s3 = $p.Slides.Item(3)

# Ovals nested inside the "Group 34" group (shape index 4) carrying the
# Alpha / RC labels on the two mini timelines.
$grp = $s3.Shapes.Item(4)
Replace-Substring $grp.GroupItems.Item(5) "Alpha" "alpha"   # Oval 29 (id 30)
Replace-Substring $grp.GroupItems.Item(6) "RC" "rc"         # Oval 23 (id 24)
Replace-Substring $grp.GroupItems.Item(7) "RC" "rc"         # Oval 28 (id 29)
Replace-Substring $grp.GroupItems.Item(8) "Alpha" "alpha"   # Oval 24 (id 25)

# "CAMARA API Definitions" textboxes (id 55 and id 57)
Replace-Substring $s3.Shapes.Item(9)  "Definitions" "definitions"
Replace-Substring $s3.Shapes.Item(10) "Definitions" "definitions"

# "(Operator/Vendor) Lab Implementations & Tests of Release Candidates (RC)"
# textboxes (id 58 and id 60)
Replace-Substring $s3.Shapes.Item(11) "Lab Implementations" "Lab implementations"
Replace-Substring $s3.Shapes.Item(11) "& Tests of Release " "& tests of release "
Replace-Substring $s3.Shapes.Item(11) "Candidates (RC)" "candidates (rc)"

Replace-Substring $s3.Shapes.Item(12) "Lab Implementations" "Lab implementations"
Replace-Substring $s3.Shapes.Item(12) "& Tests of Release " "& tests of release "
Replace-Substring $s3.Shapes.Item(12) "Candidates (RC)" "candidates (rc)"

# "Operator Production Deployments, Certifications & Channel Integrations"
# textboxes (id 61 and id 62)
Replace-Substring $s3.Shapes.Item(13) "Operator Production " "Operator production "
Replace-Substring $s3.Shapes.Item(13) "Deployments, Certifications " "deployments, certifications "
Replace-Substring $s3.Shapes.Item(13) "& Channel Integrations" "& channel integrations"

Replace-Substring $s3.Shapes.Item(14) "Operator Production " "Operator production "
Replace-Substring $s3.Shapes.Item(14) "Deployments, Certifications " "deployments, certifications "
Replace-Substring $s3.Shapes.Item(14) "& Channel Integrations" "& channel integrations"

# Bulleted placeholder (id 65) - release cadence wording
Replace-Substring $s3.Shapes.Item(15) "Fall Release ~ September" "Fall (in September)"
Replace-Substring $s3.Shapes.Item(15) "Spring Release ~ March" "Spring (in March)"

# ---------------------------------------------------------------
# Slide 5 - "Beyond Meta-release - Upcoming APIs" slide
# ---------------------------------------------------------------
$s5 = $p.Slides.Item(5)
Replace-Substring $s5.Shapes.Item(2) "Beyond Meta-release $([char]0x2013) " "Beyond Fall24 Meta-release $([char]0x2013) "
